# Add 5 new rows (97-101) of complaint-resolution data to the end of the
# table on the active sheet, matching the shape of the existing rows
# (columns A=Expediente, B=Materia, C=Resolucion, D=Descripcion (blank),
# E=URL).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ A = 202200087538; B = "Demora en la instalacion del suministro"; C = "Responsabilidad de la instaladora"; E = "/1RA5B6RV-QBRxaPkXbjiSzCgRUzv5Il_n/" },
    @{ A = 202200087540; B = "Cobro injustificado"; C = "Si hubo notificacion"; E = "/1Y25FosqhhnOIo93SYq70vqZnP3ryiIBb/" },
    @{ A = 202200087557; B = "Demora en reconexion"; C = "Visita infructuosa"; E = "/10ooyDMqFY2BMSDJYuVNueLGRmkYdPCUb/" },
    @{ A = 202200087789; B = "Demora en reconexion"; C = "Visita infructuosa"; E = "/1z29aD4oPvguGGHDJcyGNq_Myq4cGnBVg/" },
    @{ A = 202200088517; B = "Demora en reconexion"; C = "Hubieron 2 recibos pendientes de pago"; E = "/1RafrnefXenrOyspGtMUDDU1tFNlz4bIn/" }
)

$startRow = 97
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $item = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 5).Value = $item.E
}
